$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from the very start of the document
#    (it currently sits around "Description SimpleSynth").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Fix the grammar slip: "calls this function" -> "call this function"
$d.Content.Find.Execute(
    "calls this function",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "call this function", 2) | Out-Null

# 3. Re-insert the _GoBack bookmark further down, splitting the run
#    "When a key is pressed this is forwarded to the SimpleSynth (onKeyEvent)"
#    right after "this i" (i.e. before "s forwarded...").
$target = $d.Content
$target.Find.Execute(
    "When a key is pressed this is forwarded to the SimpleSynth (onKeyEvent)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "", 0) | Out-Null

$splitPoint = $target.Start + "When a key is pressed this i".Length
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
